# Cryptos list update (Sun Jan 28 07:32:13 UTC 2024) applied via Excel COM interop.
# For the Price column (D), values that look like plain numbers (e.g. "157.26")
# are written with a leading apostrophe so Excel stores them as literal text
# (matching the workbook's existing inlineStr/text-typed Price cells) instead of
# silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.477.48'
$ws.Range("E2").Value = '  +2.04%  '

$ws.Range("D3").Value = '2.288.60'
$ws.Range("E3").Value = '  +1.24%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''157.26'
$ws.Range("E5").Value = '  +15,616.25%  '

$ws.Range("E6").Value = '  +1.46%  '

$ws.Range("D7").Value = '''95.98'
$ws.Range("E7").Value = '  +5.50%  '

$ws.Range("E8").Value = '  +0.66%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = '''0.496'
$ws.Range("E10").Value = '  +3.63%  '

$ws.Range("D11").Value = '''36.20'
$ws.Range("E11").Value = '  +13.03%  '

$ws.Range("D12").Value = '''0.0805'
$ws.Range("E12").Value = '  +1.33%  '

$ws.Range("E13").Value = '  -1.91%  '

$ws.Range("E14").Value = '  +2.81%  '

$ws.Range("D15").Value = '2.645.70'
$ws.Range("E15").Value = '  +1.34%  '

$ws.Range("D16").Value = '''14.53'
$ws.Range("E16").Value = '  +2.81%  '

$ws.Range("D17").Value = '2.307.44'
$ws.Range("E17").Value = '  +1.44%  '

$ws.Range("D18").Value = '''0.806'
$ws.Range("E18").Value = '  +6.47%  '

$ws.Range("D19").Value = '42.402.86'
$ws.Range("E19").Value = '  +2.01%  '

$ws.Range("D20").Value = '''12.67'
$ws.Range("E20").Value = '  +2.91%  '

$ws.Range("D21").Value = '0.0₃0920'
$ws.Range("E21").Value = '  +2.24%  '

$ws.Range("E22").Value = '  +2.38%  '

$ws.Range("D23").Value = '''68.12'
$ws.Range("E23").Value = '  +2.48%  '

$ws.Range("D24").Value = '''243.53'
$ws.Range("E24").Value = '  +1.58%  '

$ws.Range("E25").Value = '  +1.94%  '

$ws.Range("D26").Value = '''1.94'
$ws.Range("E26").Value = '  +2.24%  '

$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("D28").Value = '''24.11'
$ws.Range("E28").Value = '  +0.50%  '

$ws.Range("D29").Value = '''36.02'
$ws.Range("E29").Value = '  +5.68%  '

$ws.Range("E31").Value = '  +1.48%  '

$ws.Range("D32").Value = '''161.28'
$ws.Range("E32").Value = '  +0.38%  '

$ws.Range("D33").Value = '''5.33'
$ws.Range("E33").Value = '  +4.12%  '

$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").Value = '''0.0756'
$ws.Range("E35").Value = '  +2.17%  '

$ws.Range("E36").Value = '  +3.43%  '

$ws.Range("E37").Value = '  +5.39%  '

$ws.Range("D38").Value = '''17.30'
$ws.Range("E38").Value = '  +4.73%  '

$ws.Range("E39").Value = '  -0.28%  '

$ws.Range("E40").Value = '  +4.37%  '

$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("E42").Value = '  +7.15%  '

$ws.Range("D43").Value = '2.010.52'
$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '''2.29'
$ws.Range("E44").Value = '  +12.24%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''19.40'
$ws.Range("E45").Value = '  +0.44%  '

$ws.Range("E46").Value = '  +3.34%  '

$ws.Range("E47").Value = '  +5.94%  '

$ws.Range("D48").Value = '''10.15'
$ws.Range("E48").Value = '  +0.64%  '

$ws.Range("D49").Value = '''53.89'
$ws.Range("E49").Value = '  +5.09%  '

$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("D51").Value = '''73.02'
$ws.Range("E51").Value = '  +1.10%  '
